$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.326.50"
$ws.Range("E2").Value = "  +1.44%  "
$ws.Range("D3").Value = "2.174.83"
$ws.Range("E3").Value = "  -0.05%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "253.13"
$ws.Range("E5").Value = "  +6.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.610"
$ws.Range("E6").Value = "  -0.44%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "73.93"
$ws.Range("E7").Value = "  +1.73%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.581"
$ws.Range("E9").Value = "  +0.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.73"
$ws.Range("E10").Value = "  +0.98%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0910"
$ws.Range("E11").Value = "  -0.16%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.75"
$ws.Range("E12").Value = "  +0.86%  "
$ws.Range("E13").Value = "  +0.33%  "
$ws.Range("D14").Value = "2.502.86"
$ws.Range("E14").Value = "  -0.03%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.15"
$ws.Range("E15").Value = "  -1.68%  "
$ws.Range("D16").Value = "2.182.68"
$ws.Range("E16").Value = "  -0.19%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.766"
$ws.Range("E17").Value = "  -2.19%  "
$ws.Range("D18").Value = "42.241.86"
$ws.Range("E18").Value = "  +1.43%  "
$ws.Range("E19").Value = "  -0.75%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.53"
$ws.Range("E20").Value = "  +0.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.86"
$ws.Range("E21").Value = "  +0.74%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "226.36"
$ws.Range("E22").Value = "  -0.13%  "
$ws.Range("E23").Value = "  +5.30%  "
$ws.Range("E24").Value = "  -5.69%  "
$ws.Range("E25").Value = "  -0.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.46"
$ws.Range("E26").Value = "  -2.97%  "
$ws.Range("E27").Value = "  +1.87%  "
$ws.Range("B28").Value = "PancakeSwap"
$ws.Range("C28").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.19"
$ws.Range("E28").Value = "  -0.33%  "
$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "37.01"
$ws.Range("E29").Value = "  +12.54%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.14"
$ws.Range("E30").Value = "  -1.83%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "169.08"
$ws.Range("E31").Value = "  -0.99%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.98"
$ws.Range("E32").Value = "  +0.46%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0805"
$ws.Range("E33").Value = "  +3.54%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.11"
$ws.Range("E34").Value = "  -3.82%  "
$ws.Range("E35").Value = "  -0.41%  "
$ws.Range("E36").Value = "  +4.01%  "
$ws.Range("E37").Value = "  -2.19%  "
$ws.Range("E38").Value = "  +8.04%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "11.93"
$ws.Range("E39").Value = "  -1.45%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.05"
$ws.Range("E40").Value = "  -2.21%  "
$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.195"
$ws.Range("E41").Value = "  +2.61%  "
$ws.Range("B42").Value = "MultiversX"
$ws.Range("C42").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "59.40"
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("E43").Value = "  -3.98%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "102.45"
$ws.Range("E44").Value = "  +5.15%  "
$ws.Range("B45").Value = "WOONetwork"
$ws.Range("C45").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.467"
$ws.Range("E45").Value = "  +15.38%  "
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.45"
$ws.Range("E46").Value = "  +11.13%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.26"
$ws.Range("E47").Value = "  -2.55%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0970"
$ws.Range("E48").Value = "  +0.46%  "
$ws.Range("E49").Value = "  +0.45%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.13"
$ws.Range("E50").Value = "  +0.90%  "
$ws.Range("E51").Value = "  +0.68%  "
